# Add new columns I ("I0") and J ("IF") to Sheet1, mirroring the header
# style already used by the other header cells (B1:H1), and fill in the
# per-row numeric values for rows 2-44.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---------------------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header formatting (bold / border / centered) used by the
# existing header cells, e.g. H1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# --- Data rows (2-44) ----------------------------------------------------
$iValues = @(9,9,8,7,8,6,7,8,6,7,8,8,6,8,5,9,9,5,6,8,7,7,6,6,6,5,4,4,9,8,5,6,8,3,4,8,8,8,8,4,4,5,3)
$jValues = @(9,9,8,8,9,8,7,8,6,7,8,8,6,8,6,9,9,6,6,8,7,8,6,6,6,5,4,4,9,8,5,6,8,4,5,8,8,8,8,4,4,5,3)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
